$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50, pushing existing rows 50-78 down to 51-79
$ws.Rows.Item(50).Insert()

# Fill the new row 50 with the new weekly record
$ws.Range("A50").Value = 8
$ws.Range("B50").Value = "Terminal La Palmera de La Serena"
$ws.Range("C50").Value = "Coquimbo"
$ws.Range("D50").Value = (Get-Date -Year 2021 -Month 9 -Day 13 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("E50").Value = 4
$ws.Range("F50").Value = 100112001
$ws.Range("G50").Value = "Berenjena"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 560
$ws.Range("K50").Value = 9000
$ws.Range("L50").Value = 10000
$ws.Range("M50").Value = 9500
$ws.Range("N50").Value = "`$/caja 60 unidades"
$ws.Range("O50").Value = "Región de Arica y Parinacota"
$ws.Range("P50").Value = 158
$ws.Range("Q50").Value = 60
$ws.Range("R50").Value = "Hortaliza"
